$d = $word.ActiveDocument

# Locate the end of the sentence we need to append after ("...telephone.").
$r = $d.Content
$r.Find.Execute("telephone.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$pos = $r.Start

# --- Two manual line breaks (plain formatting, merges with preceding text) ---
$ins = $d.Range($pos, $pos)
$ins.InsertAfter([char]11)
$pos = $pos + 1

$ins = $d.Range($pos, $pos)
$ins.InsertAfter([char]11)
$pos = $pos + 1

# --- Insert the whole new log-entry as plain text first, then layer bold /
#     underline on top of the relevant leading substrings so that we never
#     have to "turn formatting back off" (which would stamp explicit
#     <w:u w:val="none"/> / <w:b w:val="0"/> markers not present in the
#     target). ---
$dateText = "26/03/23: "
$timeText = "(16:12 PM ) "
$bodyText = "Working on insertAiport, almost had the same approach like the insertCustomers, I just added to the already made readField to validate and check the airportNumber, airportName and airportLocation."
$fullText = $dateText + $timeText + $bodyText

$ins = $d.Range($pos, $pos)
$ins.InsertAfter($fullText)

$dateStart = $pos
$dateEnd = $dateStart + $dateText.Length
$seg = $d.Range($dateStart, $dateEnd)
$seg.Font.Bold = 1
$seg.Font.Underline = 1

$timeStart = $dateEnd
$timeEnd = $timeStart + $timeText.Length
$seg = $d.Range($timeStart, $timeEnd)
$seg.Font.Bold = 1

$pos = $timeEnd + $bodyText.Length

Write-Output "done"
